$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (e.g. "2.60" -> 2.6).
foreach ($addr in @("D4", "D5", "D6", "D7", "D10", "D11", "D13", "D16", "D17", "D19", "D21", "D22", "D23", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D46", "D47", "D48", "D49", "D50")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated coin data scraped on Thu Jan 18 22:08:03 UTC 2024
$ws.Range("D2").Value = '41.275.22'
$ws.Range("E2").Value = '  -3.30%  '
$ws.Range("D3").Value = '2.459.65'
$ws.Range("E3").Value = '  -2.65%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '312.27'
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("D6").Value = '94.21'
$ws.Range("E6").Value = '  -7.42%  '
$ws.Range("D7").Value = '0.551'
$ws.Range("E7").Value = '  -3.41%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -4.68%  '
$ws.Range("D10").Value = '33.45'
$ws.Range("E10").Value = '  -7.83%  '
$ws.Range("D11").Value = '0.0781'
$ws.Range("E11").Value = '  -2.71%  '
$ws.Range("E12").Value = '  -0.19%  '
$ws.Range("D13").Value = '6.97'
$ws.Range("E13").Value = '  -5.05%  '
$ws.Range("D14").Value = '2.837.77'
$ws.Range("E14").Value = '  -2.76%  '
$ws.Range("D15").Value = '2.472.44'
$ws.Range("E15").Value = '  -0.82%  '
$ws.Range("D16").Value = '14.62'
$ws.Range("E16").Value = '  -6.48%  '
$ws.Range("D17").Value = '0.787'
$ws.Range("E17").Value = '  -2.74%  '
$ws.Range("D18").Value = '41.207.39'
$ws.Range("E18").Value = '  -3.39%  '
$ws.Range("D19").Value = '6.35'
$ws.Range("E20").Value = '  -3.40%  '
$ws.Range("D21").Value = '11.52'
$ws.Range("E21").Value = '  -5.65%  '
$ws.Range("D22").Value = '67.60'
$ws.Range("E22").Value = '  -2.64%  '
$ws.Range("D23").Value = '236.95'
$ws.Range("E23").Value = '  -3.07%  '
$ws.Range("E24").Value = '  -4.09%  '
$ws.Range("E25").Value = '  -5.53%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").Value = '24.47'
$ws.Range("E27").Value = '  -6.03%  '
$ws.Range("E28").Value = '  -4.68%  '
$ws.Range("D29").Value = '9.70'
$ws.Range("E29").Value = '  -4.57%  '
$ws.Range("D30").Value = '36.31'
$ws.Range("E30").Value = '  -6.98%  '
$ws.Range("D31").Value = '152.96'
$ws.Range("E31").Value = '  -1.78%  '
$ws.Range("D32").Value = '5.60'
$ws.Range("E32").Value = '  -3.56%  '
$ws.Range("B33").Value = 'ApeXProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D33").Value = '2.60'
$ws.Range("E33").Value = '  -5.74%  '
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '2.57'
$ws.Range("E34").Value = '  -2.05%  '
$ws.Range("D35").Value = '0.0754'
$ws.Range("E35").Value = '  -4.73%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.91'
$ws.Range("E36").Value = '  -5.96%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").Value = '3.01'
$ws.Range("E37").Value = '  -5.48%  '
$ws.Range("D38").Value = '17.08'
$ws.Range("E38").Value = '  -6.15%  '
$ws.Range("D39").Value = '0.104'
$ws.Range("E39").Value = '  -7.19%  '
$ws.Range("E40").Value = '  -3.98%  '
$ws.Range("D41").Value = '4.22'
$ws.Range("E41").Value = '  -2.13%  '
$ws.Range("D42").Value = '21.18'
$ws.Range("E42").Value = '  -4.88%  '
$ws.Range("E43").Value = '  +0.16%  '
$ws.Range("D44").Value = '1.960.23'
$ws.Range("E44").Value = '  -1.44%  '
$ws.Range("E45").Value = '  -4.86%  '
$ws.Range("D46").Value = '3.08'
$ws.Range("E46").Value = '  -7.14%  '
$ws.Range("D47").Value = '8.72'
$ws.Range("E47").Value = '  -1.59%  '
$ws.Range("D48").Value = '70.05'
$ws.Range("E48").Value = '  -3.65%  '
$ws.Range("D49").Value = '76.30'
$ws.Range("E49").Value = '  -5.32%  '
$ws.Range("D50").Value = '97.52'
$ws.Range("E50").Value = '  -3.41%  '
$ws.Range("E51").Value = '  -5.94%  '
